$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 23.79854181694116
$ws.Cells.Item(2, 3).Value = 10.65238226175392
$ws.Cells.Item(2, 4).Value = 12.43105119600247
$ws.Cells.Item(2, 5).Value = 11.34510114690246
$ws.Cells.Item(2, 7).Value = 3.852433114933596
$ws.Cells.Item(2, 9).Value = 45.48020209702018
$ws.Cells.Item(2, 10).Value = 7.184983035267148
$ws.Cells.Item(2, 11).Value = 21.75092396705796
$ws.Cells.Item(2, 12).Value = 14.1892764718536
$ws.Cells.Item(2, 13).Value = 21.87761237303273
$ws.Cells.Item(2, 14).Value = 29.35031405693727
$ws.Cells.Item(3, 2).Value = 23.73877775117528
$ws.Cells.Item(3, 3).Value = 10.55425637808902
$ws.Cells.Item(3, 4).Value = 12.43953632187688
$ws.Cells.Item(3, 5).Value = 11.36402219750183
$ws.Cells.Item(3, 7).Value = 3.855924324268384
$ws.Cells.Item(3, 9).Value = 45.46897988540396
$ws.Cells.Item(3, 10).Value = 7.176801046840313
$ws.Cells.Item(3, 11).Value = 21.71043370143645
$ws.Cells.Item(3, 12).Value = 14.20944476368527
$ws.Cells.Item(3, 13).Value = 21.89564302319829
$ws.Cells.Item(3, 14).Value = 29.36634804482702
$ws.Cells.Item(4, 2).Value = 23.70754924347911
$ws.Cells.Item(4, 3).Value = 10.49637386890186
$ws.Cells.Item(4, 4).Value = 12.44673780338099
$ws.Cells.Item(4, 5).Value = 11.3765111352738
$ws.Cells.Item(4, 7).Value = 3.858180330190594
$ws.Cells.Item(4, 9).Value = 45.46765934141258
$ws.Cells.Item(4, 10).Value = 7.171642133963439
$ws.Cells.Item(4, 11).Value = 21.69028350382921
$ws.Cells.Item(4, 12).Value = 14.22365102108455
$ws.Cells.Item(4, 13).Value = 21.91030613431866
$ws.Cells.Item(4, 14).Value = 29.37792904039561
$ws.Cells.Item(5, 2).Value = 23.69620690623039
$ws.Cells.Item(5, 3).Value = 10.47340423782687
$ws.Cells.Item(5, 4).Value = 12.4501735016592
$ws.Cells.Item(5, 5).Value = 11.38182009121422
$ws.Cells.Item(5, 7).Value = 3.859128034222175
$ws.Cells.Item(5, 9).Value = 45.4685230619642
$ws.Cells.Item(5, 10).Value = 7.169505733651531
$ws.Cells.Item(5, 11).Value = 21.68326214155121
$ws.Cells.Item(5, 12).Value = 14.22989875531688
$ws.Cells.Item(5, 13).Value = 21.91718506861973
$ws.Cells.Item(5, 14).Value = 29.3830850631808
$ws.Cells.Item(6, 2).Value = 23.69440731307435
$ws.Cells.Item(6, 3).Value = 10.46962810943754
$ws.Cells.Item(6, 4).Value = 12.45077426249323
$ws.Cells.Item(6, 5).Value = 11.38271491829601
$ws.Cells.Item(6, 7).Value = 3.85928711582843
$ws.Cells.Item(6, 9).Value = 45.46875117819651
$ws.Cells.Item(6, 10).Value = 7.169148922837511
$ws.Cells.Item(6, 11).Value = 21.68216826095267
$ws.Cells.Item(6, 12).Value = 14.23096388833176
$ws.Cells.Item(6, 13).Value = 21.9183818900343
$ws.Cells.Item(6, 14).Value = 29.38396759410629
$ws.Cells.Item(7, 2).Value = 23.70739066427235
$ws.Cells.Item(7, 3).Value = 10.49606156129536
$ws.Cells.Item(7, 4).Value = 12.44678210951046
$ws.Cells.Item(7, 5).Value = 11.37658184382837
$ws.Cells.Item(7, 7).Value = 3.858192996282592
$ws.Cells.Item(7, 9).Value = 45.46766531261277
$ws.Cells.Item(7, 10).Value = 7.171613459952543
$ws.Cells.Item(7, 11).Value = 21.6901839866585
$ws.Cells.Item(7, 12).Value = 14.22373342325745
$ws.Cells.Item(7, 13).Value = 21.91039524719124
$ws.Cells.Item(7, 14).Value = 29.37799680816933
$ws.Cells.Item(8, 2).Value = 23.77680623282883
$ws.Cells.Item(8, 3).Value = 10.61807116539671
$ws.Cells.Item(8, 4).Value = 12.43356369227576
$ws.Cells.Item(8, 5).Value = 11.35144458429976
$ws.Cells.Item(8, 7).Value = 3.85361361908532
$ws.Cells.Item(8, 9).Value = 45.47517750080469
$ws.Cells.Item(8, 10).Value = 7.182189860699115
$ws.Cells.Item(8, 11).Value = 21.73598868590526
$ws.Cells.Item(8, 12).Value = 14.19585229755743
$ws.Cells.Item(8, 13).Value = 21.88308401910239
$ws.Cells.Item(8, 14).Value = 29.3554823341564
$ws.Cells.Item(9, 2).Value = 23.95585841484185
$ws.Cells.Item(9, 3).Value = 10.87505735529182
$ws.Cells.Item(9, 4).Value = 12.42343057121172
$ws.Cells.Item(9, 5).Value = 11.30904159846387
$ws.Cells.Item(9, 7).Value = 3.845520578257439
$ws.Cells.Item(9, 9).Value = 45.53403804889994
$ws.Cells.Item(9, 10).Value = 7.201863786065813
$ws.Cells.Item(9, 11).Value = 21.86291210451408
$ws.Cells.Item(9, 12).Value = 14.15563169085322
$ws.Cells.Item(9, 13).Value = 21.85800970066493
$ws.Cells.Item(9, 14).Value = 29.32510196508143
$ws.Cells.Item(10, 2).Value = 24.11289173401076
$ws.Cells.Item(10, 3).Value = 11.07318537814879
$ws.Cells.Item(10, 4).Value = 12.42558607573547
$ws.Cells.Item(10, 5).Value = 11.28205842782248
$ws.Cells.Item(10, 7).Value = 3.840108913866524
$ws.Cells.Item(10, 9).Value = 45.60408062248519
$ws.Cells.Item(10, 10).Value = 7.215677033136664
$ws.Cells.Item(10, 11).Value = 21.97834559362354
$ws.Cells.Item(10, 12).Value = 14.13488006663891
$ws.Cells.Item(10, 13).Value = 21.85691549722472
$ws.Cells.Item(10, 14).Value = 29.31117203138906
$ws.Cells.Item(11, 2).Value = 24.18969022732574
$ws.Cells.Item(11, 3).Value = 11.16501102038851
$ws.Cells.Item(11, 4).Value = 12.42864329765253
$ws.Cells.Item(11, 5).Value = 11.27068214981141
$ws.Cells.Item(11, 7).Value = 3.837761636023826
$ws.Cells.Item(11, 9).Value = 45.64173228432094
$ws.Cells.Item(11, 10).Value = 7.221823730948143
$ws.Cells.Item(11, 11).Value = 22.03556307875191
$ws.Cells.Item(11, 12).Value = 14.12734650457492
$ws.Cells.Item(11, 13).Value = 21.86016788672198
$ws.Cells.Item(11, 14).Value = 29.30665573110152
$ws.Cells.Item(12, 2).Value = 24.21952638020245
$ws.Cells.Item(12, 3).Value = 11.19999839600001
$ws.Cells.Item(12, 4).Value = 12.43009864873355
$ws.Cells.Item(12, 5).Value = 11.26650294558402
$ws.Cells.Item(12, 7).Value = 3.836889143885934
$ws.Cells.Item(12, 9).Value = 45.6568187385974
$ws.Cells.Item(12, 10).Value = 7.224131752768022
$ws.Cells.Item(12, 11).Value = 22.0578954594321
$ws.Cells.Item(12, 12).Value = 14.12476748394234
$ws.Cells.Item(12, 13).Value = 21.86193720474363
$ws.Cells.Item(12, 14).Value = 29.30520713748702
$ws.Cells.Item(13, 2).Value = 24.21306736935513
$ws.Cells.Item(13, 3).Value = 11.19245413891056
$ws.Cells.Item(13, 4).Value = 12.42977199202459
$ws.Cells.Item(13, 5).Value = 11.26739729327095
$ws.Cells.Item(13, 7).Value = 3.837076324097787
$ws.Cells.Item(13, 9).Value = 45.65353281572727
$ws.Cells.Item(13, 10).Value = 7.223635549020252
$ws.Cells.Item(13, 11).Value = 22.05305638064418
$ws.Cells.Item(13, 12).Value = 14.12531075236868
$ws.Cells.Item(13, 13).Value = 21.8615322632913
$ws.Cells.Item(13, 14).Value = 29.30550748460489
$ws.Cells.Item(14, 2).Value = 24.19212985121966
$ws.Cells.Item(14, 3).Value = 11.16788531025345
$ws.Cells.Item(14, 4).Value = 12.42875707108466
$ws.Cells.Item(14, 5).Value = 11.27033574652855
$ws.Cells.Item(14, 7).Value = 3.837689528019111
$ws.Cells.Item(14, 9).Value = 45.64295687372423
$ws.Cells.Item(14, 10).Value = 7.222014005500974
$ws.Cells.Item(14, 11).Value = 22.03738711082259
$ws.Cells.Item(14, 12).Value = 14.12712884242347
$ws.Cells.Item(14, 13).Value = 21.86030268295021
$ws.Cells.Item(14, 14).Value = 29.30653131192295
$ws.Cells.Item(15, 2).Value = 24.1794027152653
$ws.Cells.Item(15, 3).Value = 11.15286333082681
$ws.Cells.Item(15, 4).Value = 12.42817413367444
$ws.Cells.Item(15, 5).Value = 11.2721523868966
$ws.Cells.Item(15, 7).Value = 3.838067262274917
$ws.Cells.Item(15, 9).Value = 45.63658658747758
$ws.Cells.Item(15, 10).Value = 7.221018213928778
$ws.Cells.Item(15, 11).Value = 22.02787552655549
$ws.Cells.Item(15, 12).Value = 14.12827811690679
$ws.Cells.Item(15, 13).Value = 21.85961950433702
$ws.Cells.Item(15, 14).Value = 29.30719250284816
$ws.Cells.Item(16, 2).Value = 24.10797926816976
$ws.Cells.Item(16, 3).Value = 11.06721602192048
$ws.Cells.Item(16, 4).Value = 12.4254279871028
$ws.Cells.Item(16, 5).Value = 11.28281993527064
$ws.Cells.Item(16, 7).Value = 3.840264610106441
$ws.Cells.Item(16, 9).Value = 45.60173616785477
$ws.Cells.Item(16, 10).Value = 7.215272583735804
$ws.Cells.Item(16, 11).Value = 21.97470004525139
$ws.Cells.Item(16, 12).Value = 14.13541073492728
$ws.Cells.Item(16, 13).Value = 21.85677828125935
$ws.Cells.Item(16, 14).Value = 29.31150380266329
$ws.Cells.Item(17, 2).Value = 24.06552523859168
$ws.Cells.Item(17, 3).Value = 11.01508811718444
$ws.Cells.Item(17, 4).Value = 12.42427455442887
$ws.Cells.Item(17, 5).Value = 11.2895939316285
$ws.Cells.Item(17, 7).Value = 3.841641873402202
$ws.Cells.Item(17, 9).Value = 45.581836665053
$ws.Cells.Item(17, 10).Value = 7.21171280650937
$ws.Cells.Item(17, 11).Value = 21.94327617813033
$ws.Cells.Item(17, 12).Value = 14.14027440694358
$ws.Cells.Item(17, 13).Value = 21.85599475995812
$ws.Cells.Item(17, 14).Value = 29.31461481626426
$ws.Cells.Item(18, 2).Value = 24.04161271312897
$ws.Cells.Item(18, 3).Value = 10.98526693176301
$ws.Cells.Item(18, 4).Value = 12.42380666460985
$ws.Cells.Item(18, 5).Value = 11.29357475362035
$ws.Cells.Item(18, 7).Value = 3.842444823274141
$ws.Cells.Item(18, 9).Value = 45.57093614247677
$ws.Cells.Item(18, 10).Value = 7.209652413951277
$ws.Cells.Item(18, 11).Value = 21.92564571361892
$ws.Cells.Item(18, 12).Value = 14.14325133251209
$ws.Cells.Item(18, 13).Value = 21.85589719072392
$ws.Cells.Item(18, 14).Value = 29.31657555351393
$ws.Cells.Item(19, 2).Value = 24.03360373167596
$ws.Cells.Item(19, 3).Value = 10.97519858969158
$ws.Cells.Item(19, 4).Value = 12.42368185557513
$ws.Cells.Item(19, 5).Value = 11.29493713541149
$ws.Cells.Item(19, 7).Value = 3.84271854348046
$ws.Cells.Item(19, 9).Value = 45.56733915132424
$ws.Cells.Item(19, 10).Value = 7.2089525811749
$ws.Cells.Item(19, 11).Value = 21.91975288100593
$ws.Cells.Item(19, 12).Value = 14.14429010293868
$ws.Cells.Item(19, 13).Value = 21.85592484242521
$ws.Cells.Item(19, 14).Value = 29.31726886173337
$ws.Cells.Item(20, 2).Value = 24.06999230229608
$ws.Cells.Item(20, 3).Value = 11.02062072545446
$ws.Cells.Item(20, 4).Value = 12.42437710886542
$ws.Cells.Item(20, 5).Value = 11.28886407586704
$ws.Cells.Item(20, 7).Value = 3.841494145924489
$ws.Cells.Item(20, 9).Value = 45.58389860725234
$ws.Cells.Item(20, 10).Value = 7.212093085502827
$ws.Cells.Item(20, 11).Value = 21.94657545949121
$ws.Cells.Item(20, 12).Value = 14.13973808883708
$ws.Cells.Item(20, 13).Value = 21.85604163405036
$ws.Cells.Item(20, 14).Value = 29.31426590914647
$ws.Cells.Item(21, 2).Value = 24.19825937951995
$ws.Cells.Item(21, 3).Value = 11.17509616718308
$ws.Cells.Item(21, 4).Value = 12.42904710842066
$ws.Cells.Item(21, 5).Value = 11.26946916110206
$ws.Cells.Item(21, 7).Value = 3.837508971620969
$ws.Cells.Item(21, 9).Value = 45.6460408298126
$ws.Cells.Item(21, 10).Value = 7.222490823238094
$ws.Cells.Item(21, 11).Value = 22.04197159280145
$ws.Cells.Item(21, 12).Value = 14.12658739831139
$ws.Cells.Item(21, 13).Value = 21.86064926095123
$ws.Cells.Item(21, 14).Value = 29.30622348976478
$ws.Cells.Item(22, 2).Value = 24.28647686306668
$ws.Cells.Item(22, 3).Value = 11.27729446424002
$ws.Cells.Item(22, 4).Value = 12.43383358152142
$ws.Cells.Item(22, 5).Value = 11.25754367908475
$ws.Cells.Item(22, 7).Value = 3.83499980854147
$ws.Cells.Item(22, 9).Value = 45.69148186731567
$ws.Cells.Item(22, 10).Value = 7.22917215364776
$ws.Cells.Item(22, 11).Value = 22.10819079055446
$ws.Cells.Item(22, 12).Value = 14.11958822860165
$ws.Cells.Item(22, 13).Value = 21.86679388210693
$ws.Cells.Item(22, 14).Value = 29.30249220281398
$ws.Cells.Item(23, 2).Value = 24.23899798736771
$ws.Cells.Item(23, 3).Value = 11.22264547259915
$ws.Cells.Item(23, 4).Value = 12.4311206180088
$ws.Cells.Item(23, 5).Value = 11.26384003905115
$ws.Cells.Item(23, 7).Value = 3.836330300569245
$ws.Cells.Item(23, 9).Value = 45.6667888148929
$ws.Cells.Item(23, 10).Value = 7.225616616580331
$ws.Cells.Item(23, 11).Value = 22.07249803922535
$ws.Cells.Item(23, 12).Value = 14.12317795834505
$ws.Cells.Item(23, 13).Value = 21.86322826983319
$ws.Cells.Item(23, 14).Value = 29.30434418951004
$ws.Cells.Item(24, 2).Value = 24.06797120007772
$ws.Cells.Item(24, 3).Value = 11.01811897056142
$ws.Cells.Item(24, 4).Value = 12.42433013581715
$ws.Cells.Item(24, 5).Value = 11.2891937746241
$ws.Cells.Item(24, 7).Value = 3.841560898800407
$ws.Cells.Item(24, 9).Value = 45.58296472089826
$ws.Cells.Item(24, 10).Value = 7.211921204457431
$ws.Cells.Item(24, 11).Value = 21.94508249688965
$ws.Cells.Item(24, 12).Value = 14.13997999525838
$ws.Cells.Item(24, 13).Value = 21.85601934295179
$ws.Cells.Item(24, 14).Value = 29.31442311369303
$ws.Cells.Item(25, 2).Value = 23.90288983805125
$ws.Cells.Item(25, 3).Value = 10.80378818713066
$ws.Cells.Item(25, 4).Value = 12.42448359142259
$ws.Cells.Item(25, 5).Value = 11.31977817390588
$ws.Cells.Item(25, 7).Value = 3.847615664554892
$ws.Cells.Item(25, 9).Value = 45.51339927424979
$ws.Cells.Item(25, 10).Value = 7.196654703526925
$ws.Cells.Item(25, 11).Value = 21.82464497929818
$ws.Cells.Item(25, 12).Value = 14.16496619053165
$ws.Cells.Item(25, 13).Value = 21.86174642690338
$ws.Cells.Item(25, 14).Value = 29.33184679790965
